$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-12 from 2023-10-08 (45207)
# to 2023-10-09 (45208), keeping existing cell formatting.
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 3).Value = 45208
}
